$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.145.36"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "1.832.66"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6650"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.60%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2932"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "1.829.75"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.981"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6677"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.086"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008378"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").Value = "29.122.04"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "226.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.79%  "

$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.164"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1411"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.610"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.511"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.108"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.042"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.190"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05320"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7590"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.869"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.674"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.53%  "

$ws.Range("D37").Value = "1.273.50"
$ws.Range("E37").Value = "  -3.10%  "

$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.726"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08716"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.954"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("D45").Value = "1.977.45"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5167"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.771"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.779"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "
